$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 14 (US S019 - "Mettre l'insert d'un jeu et ses genres sous forme de procedure stockee")
# gets a new Description comment and is marked done; row grows taller to fit the wrapped text.
$d14 = $ws.Range("D14")
$d14.Value2 = "'=> Insert mis en place avant la consigne, Finalement je met plutôt le DELETE de jeu et ses exemplaires en proc stoc"
$d14.NumberFormat = "@"
$d14.WrapText = $true

$ws.Range("E14").Value2 = $true

$ws.Rows("14:14").RowHeight = 48.75

# Move the active selection like the author left it.
$ws.Range("F14").Select()
